$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text happens to look like a plain decimal number.
# Format each as Text individually before writing so Excel keeps the literal
# string (incl. trailing zeros) instead of silently converting it to a number
# (comma-separated union Ranges only apply NumberFormat to their first area,
# so each cell is formatted one at a time rather than as a single union Range).

$ws.Range("D2").Value = "38.267.94"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.092.96"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.60"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.10"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0846"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "2.402.65"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.75"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.32"
$ws.Range("E14").Value = "  +6.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.47"
$ws.Range("E15").Value = "  +5.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.776"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "2.097.30"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "38.179.85"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.35"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.68"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.92"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +6.93%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  +9.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.73"
$ws.Range("E33").Value = "  +5.44%  "
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.41"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.53"
$ws.Range("E38").Value = "  +5.91%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.14"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").Value = "1.549.08"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.83"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0907"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.52"
$ws.Range("E48").Value = "  +4.99%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "2.287.69"
$ws.Range("E51").Value = "  +2.82%  "
